function Replace-ExactText {
    param($OldText, $NewText)
    $d = $word.ActiveDocument
    $searchStart = 0
    $count = 0
    while ($true) {
        $endPos = $d.Content.End
        if ($searchStart -ge $endPos) { break }
        $rng = $d.Range($searchStart, $endPos)
        $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $rng.Text = $NewText
        $count = $count + 1
        $searchStart = $rng.End
    }
    return $count
}

$null = Replace-ExactText "ESTO ES UN DOCUMENTO DE Prueba" "ESTO ES UN DOCUMENTO DE  Prueba"
$null = Replace-ExactText "En una pequeña aldea rodeada de montañas, vivía un joven llamado Miguel. Todos los días, subía al pico más alto para ver el amanecer. Un día, mientras estaba sentado contemplando el horizonte, una extraña figura apareció en el cielo. Era un ave enorme, más grande de lo que Miguel había visto jamás. El ave se acercó volando majestuosamente, y, al posarse a su lado, le habló en un idioma que Miguel no conocía. Sabía que algo grande estaba por suceder." "En una  pequeña aldea rodeada de montañas , vivía un joven  llamado Miguel. Todos los días , subía al pico más alto para ver el amanecer . Un día , mientras  estaba  sentado  contemplando el horizonte , una extraña figura apareció en el cielo . Era un ave enorme, más  grande de lo que Miguel había  visto  jamás . El ave se acercó  volando  majestuosamente , y, al posarse a su  lado , le  habló en un idioma que Miguel no conocía . Sabía que algo  grande  estaba por suceder ."
$null = Replace-ExactText "En una pequeña aldea rodeada de montañas, vivía un joven llamado Miguel. Todos los días, subía al pico más alto para ver el amanecer. Un día, mientras estaba sentado contemplando el horizonte, una extraña figura apareció en el cielo. Era un ave enorme, más grande de lo que Miguel había visto jamás. El ave se acercó volando majestuosamente, y, al posarse a su lado, le habló en un idioma que Miguel no conocía. Sabía que algo grande estaba por suceder." "En una  pequeña aldea rodeada de montañas , vivía un joven  llamado Miguel. Todos los días , subía al pico más alto para ver el amanecer . Un día , mientras  estaba  sentado  contemplando el horizonte , una extraña figura apareció en el cielo . Era un ave enorme, más  grande de lo que Miguel había  visto  jamás . El ave se acercó  volando  majestuosamente , y, al posarse a su  lado , le  habló en un idioma que Miguel no conocía . Sabía que algo  grande  estaba por suceder ."
$null = Replace-ExactText "En una pequeña aldea rodeada de montañas, vivía un joven llamado Miguel. Todos los días, subía al pico más alto para ver el amanecer. Un día, mientras estaba sentado contemplando el horizonte, una extraña figura apareció en el cielo. Era un ave enorme, más grande de lo que Miguel había visto jamás. El ave se acercó volando majestuosamente, y, al posarse a su lado, le habló en un idioma que Miguel no conocía. Sabía que algo grande estaba por suceder." "En una  pequeña aldea rodeada de montañas , vivía un joven  llamado Miguel. Todos los días , subía al pico más alto para ver el amanecer . Un día , mientras  estaba  sentado  contemplando el horizonte , una extraña figura apareció en el cielo . Era un ave enorme, más  grande de lo que Miguel había  visto  jamás . El ave se acercó  volando  majestuosamente , y, al posarse a su  lado , le  habló en un idioma que Miguel no conocía . Sabía que algo  grande  estaba por suceder ."
$null = Replace-ExactText "_`"Para ensamblar la mesa, siga los siguientes pasos:" "_`"Para  ensamblar la mesa, siga los siguientes  pasos :"
$null = Replace-ExactText "Desempaquete todas las piezas y verifique que estén completas." "Desempaquete todas laspiezas yverifique que esténcompletas."
$null = Replace-ExactText "Ensamble las patas de la mesa utilizando los tornillos provistos." "Ensamble laspatas de la mesautilizando lostornillos provistos."
$null = Replace-ExactText "Fije las patas al tablero principal de la mesa con los soportes incluidos." "Fije laspatas altablero principal de la mesa con lossoportes incluidos."
$null = Replace-ExactText "Asegúrese de que todos los tornillos estén bien ajustados, pero no los apriete en exceso." "Asegúrese de quetodos lostornillos esténbien ajustados ,pero no losapriete enexceso."
$null = Replace-ExactText "Finalmente, coloque la mesa en su posición y verifique que esté nivelada.`"_" "Finalmente ,coloque la mesa ensu posición yverifique queesté nivelada.`"_"
$null = Replace-ExactText "6. Texto literario (poesía)" "6.  Texto  literario ( poesía )"
$null = Replace-ExactText "`"El viento susurra entre las hojas,como un secreto que nunca se cuenta,la luna se oculta tras las nubes,y el mar, incansable, siempre sueña.Es el eco de una historia vieja,que en la noche silenciosa despierta,y aunque nadie la escucha del todo,en el alma su sombra se queda.`"" "`"El  viento  susurra entre las hojas , como un secreto que nunca se cuenta , la luna se oculta tras las nubes, y el mar, incansable, siempre  sueña . Es el eco de una historia vieja , que en la noche silenciosa despierta , y aunque  nadie la escucha del todo , en el alma  su  sombra se queda.`""
$null = Replace-ExactText "7. Texto legal (contrato básico)" "7.  Texto legal ( contrato  básico )"
$null = Replace-ExactText "`"El presente contrato de arrendamiento se celebra entre el propietario, Juan Pérez, y el arrendatario, María López. El arrendador se compromete a ceder el uso y disfrute del inmueble ubicado en Calle 123, Ciudad, por un período de 12 meses a partir del 1 de octubre de 2024. El arrendatario, por su parte, se compromete a pagar una re" "`"El  presente  contrato de arrendamiento se celebra entre el propietario , Juan Pérez, y el arrendatario , María López. El arrendador se compromete a ceder el uso y disfrute del inmueble  ubicado en Calle 123, Ciudad, por un período de 12 meses a partir del 1 de octubre de 2024. El arrendatario , por su  parte , se compromete a pagar una re"
$null = Replace-ExactText "ta mensual de `$500. Cualquier incumplimiento de las obligaciones establecidas en este contrato dará lugar a la resolución inmediata del mismo, conforme a la ley aplicable.`"" "ta mensual de `$500.  Cualquier  incumplimiento de las obligaciones  establecidas en este contrato  dará  lugar a la resolución  inmediata del mismo , conforme a la ley aplicable.`""
$null = Replace-ExactText "8. Texto académico (ensayo corto)" "8.  Texto  académico ( ensayo  corto )"
$null = Replace-ExactText "`"El impacto de la tecnología en la educación ha sido objeto de debate en los últimos años. Si bien algunos argumentan que la digitalización ha facilitado el acceso a recursos educativos y mejorado la interacción entre estudiantes y profesores, " "`"El impacto de la  tecnología en la educación ha sido  objeto de debate en los últimos  años . Si bien  algunos  argumentan que la digitalización ha facilitado el acceso a recursos educativos y mejorado la interacción entre estudiantes y profesores ,"
$null = Replace-ExactText "otros creen que el uso excesivo de dispositivos electrónicos ha generado distracciones y problemas de concentración. Este " "otros creen que el usoexcesivo dedispositivos electrónicos hagenerado distracciones yproblemas deconcentración . Este"
$null = Replace-ExactText "ensayo busca analizar ambos puntos de vista, centrándose en los estudios más recientes sobre el tema, y propondrá soluciones para maximizar los beneficios de la tecnología en el aula.`"" "ensayo buscaanalizar ambospuntos de vista,centrándose en losestudios más recientes sobre el tema, ypropondrá soluciones paramaximizar losbeneficios de latecnología en el aula.`""
$null = Replace-ExactText "9. Texto coloquial (conversación casual)" "9.  Texto  coloquial ( conversación casual)"
$null = Replace-ExactText "`"—¿Qué tal estuvo tu fin de semana?—Increíble, fuimos a la playa y el clima estuvo perfecto. ¿Y el tuyo?—Bastante tranquilo, me quedé en casa viendo películas.—Suena genial. A veces es lo que uno necesita, un buen descanso.`"" "`"—¿Qué talestuvo tufin desemana?—Increíble ,fuimos a la playa y el climaestuvo perfecto . ¿Y eltuyo?—Bastante tranquilo , mequedé en casaviendo películas.—Suena genial. A veces es lo queuno necesita , unbuen descanso.`""
$null = Replace-ExactText "10. Texto creativo (descripción " "10.  Texto  creativo ( descripción "
$null = Replace-ExactText "`"El castillo flotaba en el aire, suspendido entre nubes púrpuras y estrellas brillantes. Cada una de sus torres estaba hecha de cristal, y desde lo alto, ríos de luz descendían hacia los valles. En su interior, criaturas de fuego y hielo convivían en armonía, mientras el tiempo mismo parecía detenerse. Era un lugar donde lo imposible se volvía real.`"" "`"El  castillo  flotaba en el aire, suspendido entre nubes púrpuras y estrellas  brillantes . Cada una de sus torres estaba  hecha de cristal , y desde lo alto, ríos de luz  descendían  hacia los valles . En su interior, criaturas de fuego y hielo  convivían en armonía , mientras el tiempo  mismo  parecía  detenerse . Era un lugar  donde lo imposible se volvía real.`""
$null = Replace-ExactText "Estos textos deberían darte una variedad interesante para probar cómo se comporta tu servicio de traducción frente a diferentes géneros, estilos y niveles de formalidad." "Estos textos  deberían  darte una variedad  interesante para probar  cómo se comporta tu servicio de traducción  frente a diferentes  géneros , estilos y niveles de formalidad ."
$null = Replace-ExactText "Este sería un título pero es un elefante" "Este sería un  título pero es un elefante"
$null = Replace-ExactText "Este sería un título pero es un elefante" "Este sería un  título pero es un elefante"